$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H37").Value = 8256.5
$ws.Range("J37").Value = 8256.5
$ws.Range("L37").Value = 24769.5
$ws.Range("N37").Value = -25021.5
$ws.Range("H61").Value = 491
$ws.Range("I61").Value = 491
$ws.Range("K61").Value = 1473
$ws.Range("M61").Value = -1301
$ws.Range("H64").Value = 7046.467
$ws.Range("J64").Value = 7445.769
$ws.Range("L64").Value = 7445.769
$ws.Range("N64").Value = -7941.769
$ws.Range("H67").Value = 7046.467
$ws.Range("J67").Value = 7445.769
$ws.Range("L67").Value = 7445.769
$ws.Range("N67").Value = -9161.769
$ws.Range("H82").Value = 1991.75
$ws.Range("I82").Value = 1991.75
$ws.Range("K82").Value = 5975.25
$ws.Range("M82").Value = -5569.25
$ws.Range("H85").Value = 1991.75
$ws.Range("I85").Value = 1991.75
$ws.Range("K85").Value = 5975.25
$ws.Range("M85").Value = -4571.25
$ws.Range("H98").Value = 3219
$ws.Range("I98").Value = 2964.6191
$ws.Range("K98").Value = 2964.6191
$ws.Range("M98").Value = -1466.6191
$ws.Range("H99").Value = 1792.5
$ws.Range("J99").Value = 2708
$ws.Range("L99").Value = 8124
$ws.Range("N99").Value = -11120
$ws.Range("H101").Value = 462.22223
$ws.Range("J101").Value = 369.4
$ws.Range("L101").Value = 1108.2
$ws.Range("N101").Value = -4352.2
$ws.Range("H104").Value = 674.8333
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H107").Value = 46058.047
$ws.Range("I107").Value = 48218
$ws.Range("K107").Value = 48218
$ws.Range("M107").Value = -46298
$ws.Range("H112").Value = 2821.1333
$ws.Range("J112").Value = 2821.1333
$ws.Range("L112").Value = 8463.3999
$ws.Range("N112").Value = -10679.3999
$ws.Range("H118").Value = 6974.875
$ws.Range("I118").Value = 7741.4287
$ws.Range("K118").Value = 23224.2861
$ws.Range("M118").Value = -21567.2861
$ws.Range("H122").Value = 3219
$ws.Range("I122").Value = 2964.6191
$ws.Range("K122").Value = 8893.8573
$ws.Range("M122").Value = -6443.8573
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()
$ws.Range("H132").Value = 1715.1351
$ws.Range("I132").Value = 1934.8667
$ws.Range("J132").Value = 773.4286
$ws.Range("K132").Value = 5804.6001
$ws.Range("L132").Value = 2320.2858
$ws.Range("M132").Value = -3274.6001
$ws.Range("N132").Value = -7380.2858
$ws.Range("H137").Value = 4240.4194
$ws.Range("I137").Value = 4194.75
$ws.Range("K137").Value = 12584.25
$ws.Range("M137").Value = -10034.25
$ws.Range("H138").Value = 4773.0454
$ws.Range("I138").Value = 1752.6
$ws.Range("J138").Value = 6335.3447
$ws.Range("K138").Value = 5257.799999999999
$ws.Range("L138").Value = 19006.0341
$ws.Range("M138").Value = -117.7999999999993
$ws.Range("N138").Value = -29286.0341

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4185.8657
$ws.Range("I32").Value = 3859.5938
$ws.Range("K32").Value = 3859.5938
$ws.Range("M32").Value = -3572.5938
$ws.Range("H74").Value = 1078.3903
$ws.Range("I74").Value = 945.3125
$ws.Range("K74").Value = 945.3125
$ws.Range("M74").Value = -71.3125
$ws.Range("H77").Value = 1078.3903
$ws.Range("I77").Value = 945.3125
$ws.Range("K77").Value = 4726.5625
$ws.Range("M77").Value = -358.5625
$ws.Range("H132").Value = 2451.9565
$ws.Range("I132").Value = 2287.814
$ws.Range("K132").Value = 6863.441999999999
$ws.Range("M132").Value = -4333.441999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 57352.79
$ws.Range("I134").Value = 4977.9287
$ws.Range("K134").Value = 14933.7861
$ws.Range("M134").Value = -12398.7861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2038.8462
$ws.Range("I132").Value = 2041.75
$ws.Range("K132").Value = 6125.25
$ws.Range("M132").Value = -3595.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 28499.25
$ws.Range("J52").Value = 28499.25
$ws.Range("L52").Value = 28499.25
$ws.Range("N52").Value = -29017.25
$ws.Range("H102").Value = 1722.2727
$ws.Range("I102").Value = 826.9167
$ws.Range("J102").Value = 2796.7
$ws.Range("K102").Value = 826.9167
$ws.Range("L102").Value = 2796.7
$ws.Range("M102").Value = 795.0833
$ws.Range("N102").Value = -6040.7
$ws.Range("H112").Value = 70000
$ws.Range("J112").Value = 70000
$ws.Range("L112").Value = 70000
$ws.Range("N112").Value = -72216
$ws.Range("H132").Value = 40771.605
$ws.Range("I132").Value = 5136.4736
$ws.Range("J132").Value = 116001.336
$ws.Range("K132").Value = 15409.4208
$ws.Range("L132").Value = 348004.008
$ws.Range("M132").Value = -12879.4208
$ws.Range("N132").Value = -353064.008
$ws.Range("H136").Value = 67220.836
$ws.Range("J136").Value = 67220.836
$ws.Range("L136").Value = 201662.508
$ws.Range("N136").Value = -206762.508

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 700
$ws.Range("I22").Value = 700
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 700
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -405
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 700
$ws.Range("I27").Value = 700
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 700
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -593
$ws.Range("N27").ClearContents()
$ws.Range("H93").Value = 43479612
$ws.Range("I93").Value = 55556756
$ws.Range("K93").Value = 55556756
$ws.Range("M93").Value = -55555508
$ws.Range("H122").Value = 3900.125
$ws.Range("I122").Value = 2375.75
$ws.Range("K122").Value = 7127.25
$ws.Range("M122").Value = -4677.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 20271
$ws.Range("J69").Value = 20271
$ws.Range("L69").Value = 20271
$ws.Range("N69").Value = -21769
$ws.Range("H72").Value = 20271
$ws.Range("J72").Value = 20271
$ws.Range("L72").Value = 60813
$ws.Range("N72").Value = -68301
$ws.Range("H100").Value = 986.7778
$ws.Range("I100").Value = 1081.3334
$ws.Range("K100").Value = 2162.6668
$ws.Range("M100").Value = -1621.6668
$ws.Range("H113").Value = 174.58333
$ws.Range("I113").Value = 154.09091
$ws.Range("J113").Value = 400
$ws.Range("K113").Value = 462.27273
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 1707.72727
$ws.Range("N113").Value = -5540
$ws.Range("H132").Value = 15509.986
$ws.Range("I132").Value = 1475.9811
$ws.Range("J132").Value = 52700.1
$ws.Range("K132").Value = 4427.9433
$ws.Range("L132").Value = 158100.3
$ws.Range("M132").Value = -1897.9433
$ws.Range("N132").Value = -163160.3
$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200
